$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Anexo 1 (sheet1): add a new data row (row 11) below the existing table.
# Copy the format from the last existing row (10) first so the new row
# picks up the same cell style (border + centered alignment) as its
# neighbours, then fill in the values. Filling F (mask) before C/D/E
# reproduces the shared-string insertion order seen in the workbook.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Anexo 1")
$ws1.Range("B10:H10").Copy()
$ws1.Range("B11:H11").PasteSpecial(-4122)

$ws1.Cells.Item(11, 6).Value = "255.255.255.252"
$ws1.Cells.Item(11, 3).Value = "172.16.1.232"
$ws1.Cells.Item(11, 4).Value = "172.16.1.235"
$ws1.Cells.Item(11, 5).Value = "172.16.1.233"
$ws1.Cells.Item(11, 2).Value = 9
$ws1.Cells.Item(11, 7).Value = 0
$ws1.Cells.Item(11, 8).Value = 2

# ---------------------------------------------------------------------------
# Anexo 2 (sheet2): swap a couple of router-extremo labels between rows.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Anexo 2")
$ws2.Cells.Item(7, 5).Value = "R3"
$ws2.Cells.Item(8, 5).Value = "R4"
$ws2.Cells.Item(15, 3).Value = "R4"
$ws2.Cells.Item(15, 5).Value = "R7"

# ---------------------------------------------------------------------------
# Anexo 3 (sheet3): row 3's count changes, and a new row is inserted between
# the old rows 3 and 4, pushing the old row 4 down to row 5.
# Copy the format of row 4 into the not-yet-existing row 5 first (so the
# pushed-down row keeps its original style + row height), copy its values
# down, and only then overwrite row 4 with the new router entry.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Anexo 3")
$ws3.Cells.Item(3, 2).Value = 5

$ws3.Range("B4:E4").Copy()
$ws3.Range("B5:E5").PasteSpecial(-4122)
$ws3.Rows.Item(5).RowHeight = $ws3.Rows.Item(4).RowHeight

$oldC4 = $ws3.Cells.Item(4, 3).Value2
$oldD4 = $ws3.Cells.Item(4, 4).Value2
$oldE4 = $ws3.Cells.Item(4, 5).Value2
$ws3.Cells.Item(5, 2).Value = 19
$ws3.Cells.Item(5, 3).Value = $oldC4
$ws3.Cells.Item(5, 4).Value = $oldD4
$ws3.Cells.Item(5, 5).Value = $oldE4

$ws3.Cells.Item(4, 2).Value = 2
$ws3.Cells.Item(4, 3).Value = "Router"
$ws3.Cells.Item(4, 4).Value = 2811
$ws3.Cells.Item(4, 5).Value = "fastEthernet, ethernet y serial"

# ---------------------------------------------------------------------------
# View state: update the remembered selection on each sheet, then finish by
# selecting Anexo 3 so it ends up the active/visible tab.
# ---------------------------------------------------------------------------
$ws1.Range("E17").Select()
$ws2.Range("F10").Select()
$ws4 = $wb.Worksheets.Item("Anexo 4")
$ws4.Range("C16").Select()

$ws3.Range("C9").Select()
